# New Work organization v2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tag")

# Update row 19 (OTH) with new description / spanish label
$ws.Range("B19").Value = "Categoría para entidades que no encajan en clasificaciones estándar, como premios, conceptos abstractos, enfermedades, festivales, términos de tecnología, entre otros."
$ws.Range("C19").Value = "Otros"

# Move the selection to reflect where the author left off editing
$ws.Range("B20").Select()
